$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Add the required "Experimental" boolean value ("true", stored as plain text in the
# sheet, matching the source FHIR IG publisher output) next to its label in B7.
# A direct .Value = "true" assignment gets auto-coerced to a native Boolean cell, so
# instead build it through a formula and paste the computed result back as a static
# value - this preserves the literal text type ("true") rather than Boolean.
$cell = $ws.Range("B7")
$cell.Formula = "=""true"""
$cell.Copy()
$cell.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

# Update the Date value to reflect the new commit date
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"
